# Adding new RAAL model Production
# Updates computed GHI/production figures on the "Daily" and "Hourly" sheets.

$wb = $excel.ActiveWorkbook
$wsDaily  = $wb.Worksheets.Item("Daily")
$wsHourly = $wb.Worksheets.Item("Hourly")

# --- Daily sheet (row 2 totals) ---
$wsDaily.Range("G2").Value = 2605.68
$wsDaily.Range("H2").Value = 5832.57
$wsDaily.Range("I2").Value = 680.16
$wsDaily.Range("J2").Value = 1107.79
$wsDaily.Range("K2").Value = 139.92
$wsDaily.Range("L2").Value = 1055.7

# --- Hourly sheet (rows 9-19) ---
$wsHourly.Range("K9").Value = 1.32
$wsHourly.Range("M9").Value = 1.32

$wsHourly.Range("I10").Value = 367.31
$wsHourly.Range("K10").Value = 39.38
$wsHourly.Range("L10").Value = 1.79
$wsHourly.Range("M10").Value = 39.06

$wsHourly.Range("H11").Value = 215.17
$wsHourly.Range("I11").Value = 601.34
$wsHourly.Range("K11").Value = 109.06
$wsHourly.Range("L11").Value = 30.34
$wsHourly.Range("M11").Value = 100.82

$wsHourly.Range("I12").Value = 708.4400000000001
$wsHourly.Range("K12").Value = 171.16
$wsHourly.Range("L12").Value = 48.55
$wsHourly.Range("M12").Value = 153.76

$wsHourly.Range("H13").Value = 411.58
$wsHourly.Range("I13").Value = 760.6799999999999
$wsHourly.Range("K13").Value = 191.98
$wsHourly.Range("L13").Value = 28.07
$wsHourly.Range("M13").Value = 179.89

$wsHourly.Range("H14").Value = 442.25
$wsHourly.Range("I14").Value = 778.22
$wsHourly.Range("J14").Value = 92.53
$wsHourly.Range("K14").Value = 191.25
$wsHourly.Range("L14").Value = 18.84
$wsHourly.Range("M14").Value = 182.62

$wsHourly.Range("H15").Value = 420.97
$wsHourly.Range("I15").Value = 766.33
$wsHourly.Range("J15").Value = 90.62
$wsHourly.Range("K15").Value = 170.04
$wsHourly.Range("L15").Value = 10.79
$wsHourly.Range("M15").Value = 165.27

$wsHourly.Range("H16").Value = 350.17
$wsHourly.Range("I16").Value = 721.42
$wsHourly.Range("J16").Value = 83.8
$wsHourly.Range("K16").Value = 123.41
$wsHourly.Range("L16").Value = 1.54
$wsHourly.Range("M16").Value = 122.78

$wsHourly.Range("H17").Value = 238.7
$wsHourly.Range("I17").Value = 627.24
$wsHourly.Range("K17").Value = 69.23
$wsHourly.Range("M17").Value = 69.23

$wsHourly.Range("I18").Value = 428.06
$wsHourly.Range("J18").Value = 48.22
$wsHourly.Range("K18").Value = 37.52
$wsHourly.Range("L18").Value = 0
$wsHourly.Range("M18").Value = 37.52

$wsHourly.Range("K19").Value = 3.44
$wsHourly.Range("M19").Value = 3.44
